$d = $word.ActiveDocument

# Locate the block to rewrite by content rather than a hard-coded paragraph
# index: it runs from the start of the "+ Quên mật khẩu" paragraph through
# the end of the last "+ Quản lí ..." paragraph ("... phương thức thanh
# toán").
$rQuen = $d.Content
$rQuen.Find.ClearFormatting()
[void]$rQuen.Find.Execute("+ Quên mật khẩu")
$pQuen = $rQuen.Paragraphs(1)

$rLast = $d.Content
$rLast.Find.ClearFormatting()
[void]$rLast.Find.Execute("+ Quản lí phương thức thanh toán")
$pLastMgmt = $rLast.Paragraphs(1)

$r = $d.Range($pQuen.Range.Start, $pLastMgmt.Range.End)

# Rewriting this whole span in a single InsertXML() call lets the host's
# "reattach a surviving bookmark" behaviour drop the old (now unreferenced)
# _GoBack bookmark instead of re-anchoring it on the rewritten
# "+ Quên mật khẩu" paragraph, while the new bookmark we declare explicitly
# on "+ Xóa sản phẩm" becomes the live _GoBack.
$pPr = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:tabs><w:tab w:val="center" w:pos="5179"/></w:tabs><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>'
$pPrAdmin = '<w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr><w:tabs><w:tab w:val="center" w:pos="5179"/></w:tabs><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr>'
$rPr = '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr>'

$xml = '<?xml version="1.0" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body>' +
    '<w:p>' + $pPr + $rPr + '<w:t>+ Quên mật khẩu</w:t></w:r></w:p>' +
    '<w:p>' + $pPrAdmin + $rPr + '<w:t>Use case (Admin):</w:t></w:r></w:p>' +
    '<w:p>' + $pPr + $rPr + '<w:t>+ Thêm sản phẩm</w:t></w:r></w:p>' +
    '<w:p>' + $pPr + $rPr + '<w:t>+ Sửa sản phẩm</w:t></w:r></w:p>' +
    '<w:p>' + $pPr + $rPr + '<w:t>+ Xóa sản phẩm</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>' +
    '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

[void]$r.InsertXML($xml)
